$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 becomes text "male", B2 becomes -0.102
$ws.Range("A2").Value = "male"
$ws.Range("B2").Value = -0.102

# Row 4 (A4="male", B4=0.283) is removed entirely
$ws.Rows("4").Delete()
